# Holly fixed the "harvester" column (column B) in rnaSamples: the values
# had incorrectly carried over the "Retrofitted_1658" text (which actually
# belongs in rnaPreparer, column E). Update harvester to "S.GISH" for every
# data row, leaving rnaPreparer (E) and rnaPrepMethod (G) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row so the fix covers every data row regardless
# of how many rows the sheet actually has.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
if ($lastRow -lt 2) {
    $lastRow = 2
}

# Row 1 is the header ("harvester"); data starts on row 2.
$ws.Range("B2:B" + $lastRow).Value = "S.GISH"
